$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RR")

# Extend the recycling-rate "target" values (column F) from the EU-only
# zero baseline to 0.01 for the historical years 2000-2023 (rows 3-26),
# matching the other regions' columns.
$ws.Range("F3:F26").Value = 0.01

# Apply an integer number format to column A (year) for the full table,
# consistent across all regions now that installation data spans them.
$ws.Range("A3:A103").NumberFormat = "0"

# Make "RR" the active/selected sheet with F3:F26 selected (mirrors the
# region coverage just extended), and move selection off "region".
$ws.Activate()
$null = $ws.Range("F3:F26").Select()
